# Tarefas.xlsx update
# - Task "Floyd Warshall" (row 3) status changes from "Fase final de
#   implementação" to "Aguarda testes".
# - New task row added for the Dijkstra algorithm implementation, assigned
#   to Tiago, status "Iniciado".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status of the existing "floyd warshall" task (row 3, column E).
$ws.Range("E3").Value = "Aguarda testes"

# Append the new task row (row 10).
$ws.Range("A10").Value = "Implementação do algoritmo de Dijkstra com retorno de distâcia de vertices"
$ws.Range("B10").Value = "Tiago"
$ws.Range("E10").Value = "Iniciado"

# Move/keep the active selection on the newly added row, as in the source file.
$ws.Range("E10").Select()
